$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3236023485660553
$ws.Range("B1").Value = 0.5828288197517395
$ws.Range("C1").Value = 3.508979320526123
$ws.Range("D1").Value = 1.590747356414795
$ws.Range("E1").Value = 1.116180539131165
